# Generate Report for Handoff
# Updates the status/handoff-date info for the "842c3a25-2135-4d40-8894-63400e4d118a.md"
# row (row 3) across the Overview, zh-cn and de-de sheets to reflect that the
# file is now "Ready for handoff".

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-27-18 05:27:34"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-18 05:27:31"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-18 05:27:34"
